$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing quarterly data (previously in D:K) to F:M.
# This matches how a new quarter's data is prepended in this financial model workbook.
$ws.Columns("D:E").Insert()

# The inserted columns default to a blank style; restore the correct number formats
# (date format for header rows, number format for data rows) by copying formats from
# column F (which now holds what used to be column D, i.e. the adjacent existing data).
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns (D = quarter ending 2018-12-31, E = quarter ending 2018-09-30)
# with the newly reported financial figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 381400
$ws.Range("E8").Value = 323800
$ws.Range("D9").Value = 5800
$ws.Range("E9").Value = 2400
$ws.Range("D10").Value = 375600
$ws.Range("E10").Value = 321400
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 800
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 45600
$ws.Range("E15").Value = 51200
$ws.Range("D17").Value = 315500
$ws.Range("E17").Value = 295000
$ws.Range("D18").Value = 65900
$ws.Range("E18").Value = 28800
$ws.Range("D20").Value = 800
$ws.Range("E20").Value = 400
$ws.Range("D21").Value = 112400
$ws.Range("E21").Value = 80500
$ws.Range("D22").Value = 20400
$ws.Range("E22").Value = 20900
$ws.Range("D23").Value = 46300
$ws.Range("E23").Value = 8400
$ws.Range("D24").Value = -2500
$ws.Range("E24").Value = -700
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 48800
$ws.Range("E26").Value = 9100
$ws.Range("D27").Value = 42400
$ws.Range("E27").Value = 5700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -2800
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -800
$ws.Range("E32").Value = -400
$ws.Range("D33").Value = 39700
$ws.Range("E33").Value = 5700
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 39700
$ws.Range("E35").Value = 5700
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 235700
$ws.Range("E41").Value = 250500
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 128900
$ws.Range("E43").Value = 107900
$ws.Range("D44").Value = 48600
$ws.Range("E44").Value = 52500
$ws.Range("D45").Value = 60700
$ws.Range("E45").Value = 58100
$ws.Range("D46").Value = 474000
$ws.Range("E46").Value = 468900
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 2106500
$ws.Range("E48").Value = 2103700
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 180200
$ws.Range("E52").Value = 173300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2760700
$ws.Range("E54").Value = 2745900
$ws.Range("D57").Value = 130900
$ws.Range("E57").Value = 102400
$ws.Range("D58").Value = 134800
$ws.Range("E58").Value = 20900
$ws.Range("D59").Value = 226400
$ws.Range("E59").Value = 249800
$ws.Range("D60").Value = 492200
$ws.Range("E60").Value = 373100
$ws.Range("D61").Value = 734200
$ws.Range("E61").Value = 857000
$ws.Range("D62").Value = 982700
$ws.Range("E62").Value = 1053400
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 2350800
$ws.Range("E66").Value = 2425000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 182100
$ws.Range("E72").Value = 148600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 409900
$ws.Range("E76").Value = 320900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 39700
$ws.Range("E81").Value = 5700
$ws.Range("D83").Value = 45600
$ws.Range("E83").Value = 51200
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 83300
$ws.Range("E89").Value = 52100
$ws.Range("D91").Value = -48900
$ws.Range("E91").Value = -40700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -58200
$ws.Range("E94").Value = -40500
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -36500
$ws.Range("E100").Value = -19700
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -11400
$ws.Range("E102").Value = -8200

